# Update "想去人数" (interested-attendee count) figures that increased by
# small amounts since the last scrape, on both the "展览" sheet and the
# "全部类型" aggregate sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 563
$ws1.Range("F4").Value = 210
$ws1.Range("F6").Value = 517
$ws1.Range("F10").Value = 6908
$ws1.Range("F13").Value = 3268
$ws1.Range("F14").Value = 222
$ws1.Range("F15").Value = 398
$ws1.Range("F17").Value = 567
$ws1.Range("F18").Value = 38

# Sheet "全部类型" (all types, combined)
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F5").Value = 563
$ws2.Range("F6").Value = 210
$ws2.Range("F8").Value = 517
$ws2.Range("F13").Value = 6908
$ws2.Range("F17").Value = 3268
$ws2.Range("F18").Value = 222
$ws2.Range("F19").Value = 398
$ws2.Range("F21").Value = 567
$ws2.Range("F22").Value = 38
